# ============================================================================
# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# Target layout (tab order):
#   1. Player Info        (new)
#   2. ODI Batting         (existing, renamed column + new link->code values)
#   3. ODI Bowling         (existing, renamed column + new link->code values)
#   4. ODI Batting Extra   (new)
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# 1. Insert the new "Player Info" sheet before everything else.
#    NOTE: worksheet object references in this host re-resolve by *position*,
#    not identity, so every sheet handle used below is (re-)fetched by name
#    *after* all Add() calls that could shift indices have completed.
# ----------------------------------------------------------------------------
$playerInfoWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value)
$playerInfoWs.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfoWs.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$piData = @("3911", "Douglas A J Bracewell", "Right Handed", "Right Arm Medium")
for ($c = 1; $c -le $piData.Length; $c++) {
    $cell = $playerInfoWs.Cells.Item(2, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $piData[$c - 1]
}

# ----------------------------------------------------------------------------
# 2. Append the new "ODI Batting Extra" sheet after "ODI Bowling".
# ----------------------------------------------------------------------------
$bowlingWsForInsert = $wb.Worksheets.Item("ODI Bowling")
$extraWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $bowlingWsForInsert)
$extraWs.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $exHeaders.Length; $c++) {
    $cell = $extraWs.Cells.Item(1, $c)
    $cell.Value = $exHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Row data: MATCH_CODE | BATTING_POSITION | NUM_4 | NUM_6 | PERCENT_RUNS_OF_TOTAL | MAN_OF_MATCH
# "N:" prefix -> numeric cell, "S:" prefix -> text cell (blank after the colon -> blank text cell).
$exRows = @(
    @("S:3342", "S:",  "S:",  "S:",  "S:",      "S:NO"),
    @("S:3366", "N:9", "S:0", "S:0", "S:3.23%", "S:NO"),
    @("S:3370", "S:",  "S:",  "S:",  "S:",      "S:NO"),
    @("S:3385", "S:",  "S:",  "S:",  "S:",      "S:NO"),
    @("S:3429", "N:9", "S:0", "S:0", "S:2.08%", "S:NO"),
    @("S:3503", "N:9", "S:",  "S:",  "S:",      "S:NO"),
    @("S:3827", "N:8", "S:",  "S:",  "S:",      "S:NO"),
    @("S:3828", "N:8", "S:1", "S:0", "S:5.88%", "S:NO"),
    @("S:3865", "N:8", "S:",  "S:",  "S:",      "S:NO"),
    @("S:3866", "N:8", "S:",  "S:",  "S:",      "S:NO"),
    @("S:3868", "S:",  "S:",  "S:",  "S:",      "S:NO"),
    @("S:3888", "N:8", "S:0", "S:0", "S:0.81%", "S:NO"),
    @("S:3951", "S:",  "S:",  "S:",  "S:",      "S:NO"),
    @("S:4100", "N:8", "S:",  "S:",  "S:",      "S:YES"),
    @("S:4101", "N:8", "S:1", "S:0", "S:1.54%", "S:NO"),
    @("S:4239", "N:8", "S:1", "S:0", "S:4.46%", "S:NO"),
    @("S:4242", "S:",  "S:",  "S:",  "S:",      "S:NO"),
    @("S:4245", "N:8", "S:1", "S:0", "S:6.17%", "S:NO"),
    @("S:4566", "N:8", "S:2", "S:0", "S:15.53%","S:NO"),
    @("S:4568", "S:",  "S:",  "S:",  "S:",      "S:NO")
)

$r = 2
foreach ($rowVals in $exRows) {
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $raw = $rowVals[$c - 1]
        $kind = $raw.Substring(0, 1)
        $val = $raw.Substring(2)
        $cell = $extraWs.Cells.Item($r, $c)
        if ($kind -eq "N") {
            $cell.Value = [double]$val
        } else {
            $cell.NumberFormat = "@"
            $cell.Value = $val
        }
    }
    $r = $r + 1
}

# ----------------------------------------------------------------------------
# 3. "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE, replace the link
#    values with the bare numeric match code, and drop the previously-blank
#    INNING_NUMBER cells entirely (did-not-bat rows).
# ----------------------------------------------------------------------------
$battingWs = $wb.Worksheets.Item("ODI Batting")
$battingWs.Range("D1").Value = "MATCH_CODE"

$matchCodes = @("3340","3342","3366","3370","3385","3429","3503","3827","3828","3865","3866","3868","3888","3951","4100","4101","4239","4242","4245","4566","4568")
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $cell = $battingWs.Cells.Item(2 + $i, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}

$blankInningRows = @(2, 3, 8, 9, 11, 12, 16)
foreach ($row in $blankInningRows) {
    $battingWs.Cells.Item($row, 2).ClearContents()
}

# ----------------------------------------------------------------------------
# 4. "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE and replace the link
#    values with the bare numeric match code.
# ----------------------------------------------------------------------------
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")
$bowlingWs.Range("B1").Value = "MATCH_CODE"

for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $cell = $bowlingWs.Cells.Item(2 + $i, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}

# ----------------------------------------------------------------------------
# 5. Restore the first sheet as the active one (matches original activeTab=0).
# ----------------------------------------------------------------------------
$playerInfoWs.Activate()

Write-Output "Edit complete."
